$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing PriceChange / UpDown values for row 10
$ws.Cells.Item(10, 24).Value = 0.6499990000000011
$ws.Cells.Item(10, 25).Value = "Up"

# Add a new trade row (row 11)
$ws.Cells.Item(11, 1).Value = 42654.883275462962
$ws.Cells.Item(11, 1).NumberFormat = "m/d/yy h:mm"
$ws.Cells.Item(11, 2).Value = 22
$ws.Cells.Item(11, 3).Value = "Strong Buy"
$ws.Cells.Item(11, 4).Value = 22
$ws.Cells.Item(11, 5).Value = 8618
$ws.Cells.Item(11, 6).Value = 1244
$ws.Cells.Item(11, 7).Value = 63
$ws.Cells.Item(11, 8).Value = 36
$ws.Cells.Item(11, 9).Value = 88
$ws.Cells.Item(11, 10).Value = 11
$ws.Cells.Item(11, 11).Value = 16730
$ws.Cells.Item(11, 12).Value = 199
$ws.Cells.Item(11, 13).Value = 113
$ws.Cells.Item(11, 14).Value = 37
$ws.Cells.Item(11, 15).Value = 5
$ws.Cells.Item(11, 16).Value = "Noun"
$ws.Cells.Item(11, 17).Value = 58.438771163779279
$ws.Cells.Item(11, 18).Value = 0.49
$ws.Cells.Item(11, 19).Value = 0.0933
$ws.Cells.Item(11, 19).NumberFormat = "0.00%"
$ws.Cells.Item(11, 20).Value = 0.0249
$ws.Cells.Item(11, 20).NumberFormat = "0.00%"
$ws.Cells.Item(11, 21).Value = 2.34
$ws.Cells.Item(11, 22).Value = "N/A"
$ws.Cells.Item(11, 23).Value = 2

# Column C width adjustment (Verdict column now needs to fit "Strong Buy")
$ws.Columns.Item(3).ColumnWidth = 8.29
